$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.929189470362074
$ws.Range("D2").Value = 0.3628747910381653

# Row 3
$ws.Range("C3").Value = 0.8595303070699916
$ws.Range("D3").Value = 0.3993174793716168

# Row 4
$ws.Range("C4").Value = -0.03430736342932208
$ws.Range("D4").Value = 0.9729414277102797

# Row 5
$ws.Range("C5").Value = -1.750982909030123
$ws.Range("D5").Value = 0.09388194302898722

# Row 6
$ws.Range("C6").Value = -0.3902222016006309
$ws.Range("D6").Value = 0.700124861785838

# Row 7
$ws.Range("C7").Value = -1.511851389410856
$ws.Range("D7").Value = 0.1448041941464893

# Row 8
$ws.Range("C8").Value = -2.171763409615829
$ws.Range("D8").Value = 0.04093076912766969

# Row 9
$ws.Range("C9").Value = -1.167390996076714
$ws.Range("D9").Value = 0.2555498672096044

# Row 10
$ws.Range("C10").Value = -1.831204767259032
$ws.Range("D10").Value = 0.08064776561340259
$ws.Range("G10").Value = "No"

# Row 11
$ws.Range("C11").Value = -1.339052375300054
$ws.Range("D11").Value = 0.1942258479097894
